# Auto-generated edit script: updates market-price / profit columns (H-N)
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1877.8572
$ws.Range("I129").Value = 1539.2
$ws.Range("J129").Value = 2724.5
$ws.Range("K129").Value = 4617.6
$ws.Range("L129").Value = 8173.5
$ws.Range("M129").Value = 382.3999999999996
$ws.Range("N129").Value = -18173.5
$ws.Range("H131").Value = 9312.299999999999
$ws.Range("I131").Value = 1549.6
$ws.Range("K131").Value = 4648.799999999999
$ws.Range("M131").Value = 391.2000000000007
$ws.Range("H137").Value = 2436.95
$ws.Range("I137").Value = 2023.875
$ws.Range("K137").Value = 6071.625
$ws.Range("M137").Value = -3521.625
$ws.Range("H141").Value = 2858.6453
$ws.Range("I141").Value = 1578.4
$ws.Range("J141").Value = 5186.364
$ws.Range("K141").Value = 4735.200000000001
$ws.Range("L141").Value = 15559.092
$ws.Range("M141").Value = 444.7999999999993
$ws.Range("N141").Value = -25919.092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 19999.5
$ws.Range("J62").Value = 19999.5
$ws.Range("L62").Value = 19999.5
$ws.Range("N62").Value = -21247.5
$ws.Range("H63").Value = 2167
$ws.Range("I63").Value = 2850.6667
$ws.Range("J63").Value = 116
$ws.Range("K63").Value = 2850.6667
$ws.Range("L63").Value = 116
$ws.Range("M63").Value = -2164.6667
$ws.Range("N63").Value = -1488
$ws.Range("H65").Value = 19999.5
$ws.Range("J65").Value = 19999.5
$ws.Range("L65").Value = 59998.5
$ws.Range("N65").Value = -66238.5
$ws.Range("H66").Value = 2167
$ws.Range("I66").Value = 2850.6667
$ws.Range("J66").Value = 116
$ws.Range("K66").Value = 14253.3335
$ws.Range("L66").Value = 580
$ws.Range("M66").Value = -10821.3335
$ws.Range("N66").Value = -7444
$ws.Range("H74").Value = 2758.1667
$ws.Range("I74").Value = 2758.1667
$ws.Range("K74").Value = 2758.1667
$ws.Range("M74").Value = -1884.1667
$ws.Range("H77").Value = 2758.1667
$ws.Range("I77").Value = 2758.1667
$ws.Range("K77").Value = 13790.8335
$ws.Range("M77").Value = -9422.833500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1612.6428
$ws.Range("I86").Value = 1377.5555
$ws.Range("J86").Value = 2035.8
$ws.Range("K86").Value = 1377.5555
$ws.Range("L86").Value = 2035.8
$ws.Range("M86").Value = -254.5554999999999
$ws.Range("N86").Value = -4281.8
$ws.Range("H89").Value = 1612.6428
$ws.Range("I89").Value = 1377.5555
$ws.Range("J89").Value = 2035.8
$ws.Range("K89").Value = 6887.7775
$ws.Range("L89").Value = 10179
$ws.Range("M89").Value = -1271.7775
$ws.Range("N89").Value = -21411
$ws.Range("H107").Value = 1953.2941
$ws.Range("I107").Value = 1914.4
$ws.Range("K107").Value = 1914.4
$ws.Range("M107").Value = 5.599999999999909
$ws.Range("H134").Value = 3897.4092
$ws.Range("I134").Value = 3207.5264
$ws.Range("J134").Value = 8266.666999999999
$ws.Range("K134").Value = 9622.5792
$ws.Range("L134").Value = 24800.001
$ws.Range("M134").Value = -7087.5792
$ws.Range("N134").Value = -29870.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15444.1
$ws.Range("J31").Value = 16304.223
$ws.Range("L31").Value = 16304.223
$ws.Range("N31").Value = -16894.223
$ws.Range("H34").Value = 15444.1
$ws.Range("J34").Value = 16304.223
$ws.Range("L34").Value = 16304.223
$ws.Range("N34").Value = -16708.223
$ws.Range("H39").Value = 29349.5
$ws.Range("I39").Value = 18700
$ws.Range("K39").Value = 18700
$ws.Range("M39").Value = -18309
$ws.Range("H49").Value = 29349.5
$ws.Range("I49").Value = 18700
$ws.Range("K49").Value = 18700
$ws.Range("M49").Value = -18518
$ws.Range("H58").Value = 43279.6
$ws.Range("I58").Value = 50974.5
$ws.Range("J58").Value = 12500
$ws.Range("K58").Value = 50974.5
$ws.Range("L58").Value = 12500
$ws.Range("M58").Value = -50771.5
$ws.Range("N58").Value = -12906
$ws.Range("H107").Value = 908712.2
$ws.Range("I107").Value = 2717917
$ws.Range("J107").Value = 4109.75
$ws.Range("K107").Value = 2717917
$ws.Range("L107").Value = 4109.75
$ws.Range("M107").Value = -2715997
$ws.Range("N107").Value = -7949.75
$ws.Range("H136").Value = 43279.6
$ws.Range("I136").Value = 50974.5
$ws.Range("J136").Value = 12500
$ws.Range("K136").Value = 152923.5
$ws.Range("L136").Value = 37500
$ws.Range("M136").Value = -150373.5
$ws.Range("N136").Value = -42600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 84234.586
$ws.Range("J5").Value = 740
$ws.Range("L5").Value = 2220
$ws.Range("N5").Value = -2444
$ws.Range("H69").Value = 4682.25
$ws.Range("J69").Value = 5993
$ws.Range("L69").Value = 17979
$ws.Range("N69").Value = -19601
$ws.Range("H72").Value = 4682.25
$ws.Range("J72").Value = 5993
$ws.Range("L72").Value = 53937
$ws.Range("N72").Value = -62049
$ws.Range("H92").Value = 382.5
$ws.Range("I92").Value = 395
$ws.Range("K92").Value = 1185
$ws.Range("M92").Value = 63
$ws.Range("H108").Value = 1498.3636
$ws.Range("I108").Value = 1348.3
$ws.Range("K108").Value = 4044.9
$ws.Range("M108").Value = -1164.9
$ws.Range("H129").Value = 977.625
$ws.Range("J129").Value = 1199.8
$ws.Range("L129").Value = 3599.4
$ws.Range("N129").Value = -13599.4
$ws.Range("H135").Value = 84234.586
$ws.Range("J135").Value = 740
$ws.Range("L135").Value = 6660
$ws.Range("N135").Value = -11730
$ws.Range("H139").Value = 2069.2666
$ws.Range("I139").Value = 1287.6154
$ws.Range("K139").Value = 3862.8462
$ws.Range("M139").Value = 1277.1538
$ws.Range("H140").Value = 2922.75
$ws.Range("I140").Value = 2845.75
$ws.Range("K140").Value = 8537.25
$ws.Range("M140").Value = -3357.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 962.8125
$ws.Range("I2").Value = 1219.4546
$ws.Range("K2").Value = 1219.4546
$ws.Range("M2").Value = -1106.4546
$ws.Range("H21").Value = 27159.705
$ws.Range("I21").Value = 19999.334
$ws.Range("K21").Value = 19999.334
$ws.Range("M21").Value = -19826.334
$ws.Range("H22").Value = 33019.2
$ws.Range("I22").Value = 7550
$ws.Range("J22").Value = 49998.668
$ws.Range("K22").Value = 7550
$ws.Range("L22").Value = 49998.668
$ws.Range("M22").Value = -7021
$ws.Range("N22").Value = -51056.668
$ws.Range("H30").Value = 27159.705
$ws.Range("I30").Value = 19999.334
$ws.Range("K30").Value = 19999.334
$ws.Range("M30").Value = -19894.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9224.137000000001
$ws.Range("I40").Value = 8351
$ws.Range("J40").Value = 12192.8
$ws.Range("K40").Value = 8351
$ws.Range("L40").Value = 12192.8
$ws.Range("M40").Value = -8215
$ws.Range("N40").Value = -12464.8
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96240
$ws.Range("H68").Value = 1390654.5
$ws.Range("I68").Value = 1737818.6
$ws.Range("J68").Value = 1998
$ws.Range("K68").Value = 1737818.6
$ws.Range("L68").Value = 1998
$ws.Range("M68").Value = -1737069.6
$ws.Range("N68").Value = -3496
$ws.Range("H71").Value = 1390654.5
$ws.Range("I71").Value = 1737818.6
$ws.Range("J71").Value = 1998
$ws.Range("K71").Value = 8689093
$ws.Range("L71").Value = 9990
$ws.Range("M71").Value = -8685349
$ws.Range("N71").Value = -17478
$ws.Range("H93").Value = 1659.6666
$ws.Range("I93").Value = 1659.6666
$ws.Range("K93").Value = 1659.6666
$ws.Range("M93").Value = -411.6666
$ws.Range("H97").Value = 22999.25
$ws.Range("J97").Value = 22999.25
$ws.Range("L97").Value = 22999.25
$ws.Range("N97").Value = -24981.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 10407
$ws.Range("I29").Value = 8999
$ws.Range("K29").Value = 8999
$ws.Range("M29").Value = -8709
$ws.Range("H63").Value = 4312.25
$ws.Range("J63").Value = 4312.25
$ws.Range("L63").Value = 4312.25
$ws.Range("N63").Value = -5560.25
$ws.Range("H66").Value = 4312.25
$ws.Range("J66").Value = 4312.25
$ws.Range("L66").Value = 12936.75
$ws.Range("N66").Value = -19176.75
$ws.Range("H94").Value = 100330
$ws.Range("J94").Value = 100330
$ws.Range("L94").Value = 100330
$ws.Range("N94").Value = -102132
$ws.Range("H97").Value = 39999
$ws.Range("J97").Value = 39999
$ws.Range("L97").Value = 39999
$ws.Range("N97").Value = -41981
